# Applies the price/volume updates to the cryptos sheet (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '60.296.19'
$ws.Cells.Item(2, 5).Value = '  +3.37%  '

$ws.Cells.Item(3, 4).Value = '2.338.62'
$ws.Cells.Item(3, 5).Value = '  +1.73%  '

$ws.Cells.Item(4, 5).Value = '  -0.07%  '

$ws.Cells.Item(5, 4).Value = "'" + '545.66'
$ws.Cells.Item(5, 5).Value = '  +1.96%  '

$ws.Cells.Item(6, 4).Value = "'" + '131.46'
$ws.Cells.Item(6, 5).Value = '  -0.21%  '

$ws.Cells.Item(7, 5).Value = '  -0.07%  '

$ws.Cells.Item(8, 5).Value = '  -0.24%  '

$ws.Cells.Item(9, 4).Value = '2.334.21'
$ws.Cells.Item(9, 5).Value = '  +1.54%  '

$ws.Cells.Item(10, 5).Value = '  +1.20%  '

$ws.Cells.Item(11, 4).Value = "'" + '5.52'
$ws.Cells.Item(11, 5).Value = '  +0.72%  '

$ws.Cells.Item(12, 5).Value = '  +0.68%  '

$ws.Cells.Item(13, 5).Value = '  +0.93%  '

$ws.Cells.Item(14, 4).Value = "'" + '23.76'
$ws.Cells.Item(14, 5).Value = '  +0.90%  '

$ws.Cells.Item(15, 4).Value = '2.751.54'
$ws.Cells.Item(15, 5).Value = '  +1.59%  '

$ws.Cells.Item(16, 4).Value = '60.258.86'
$ws.Cells.Item(16, 5).Value = '  +3.44%  '

$ws.Cells.Item(17, 5).Value = '  +0.70%  '

$ws.Cells.Item(18, 4).Value = '2.338.55'
$ws.Cells.Item(18, 5).Value = '  +1.81%  '

$ws.Cells.Item(19, 4).Value = "'" + '10.62'
$ws.Cells.Item(19, 5).Value = '  +0.42%  '

$ws.Cells.Item(20, 5).Value = '  -1.32%  '

$ws.Cells.Item(21, 4).Value = "'" + '6.78'
$ws.Cells.Item(21, 5).Value = '  +5.61%  '

$ws.Cells.Item(22, 4).Value = "'" + '314.39'
$ws.Cells.Item(22, 5).Value = '  +0.50%  '

$ws.Cells.Item(23, 5).Value = '  -0.25%  '

$ws.Cells.Item(24, 4).Value = "'" + '63.73'
$ws.Cells.Item(24, 5).Value = '  +1.68%  '

$ws.Cells.Item(25, 4).Value = "'" + '0.171'
$ws.Cells.Item(25, 5).Value = '  +1.82%  '

$ws.Cells.Item(26, 5).Value = '  +0.12%  '

$ws.Cells.Item(27, 4).Value = "'" + '7.91'
$ws.Cells.Item(27, 5).Value = '  -1.55%  '

$ws.Cells.Item(28, 5).Value = '  +6.92%  '

$ws.Cells.Item(29, 5).Value = '  +1.91%  '

$ws.Cells.Item(30, 4).Value = "'" + '172.38'
$ws.Cells.Item(30, 5).Value = '  +1.08%  '

$ws.Cells.Item(31, 5).Value = '  +11.47%  '

$ws.Cells.Item(32, 4).Value = '0.0₃0731'
$ws.Cells.Item(32, 5).Value = '  +1.28%  '

$ws.Cells.Item(33, 5).Value = '  +3.16%  '

$ws.Cells.Item(34, 5).Value = '  +12.48%  '

$ws.Cells.Item(35, 5).Value = '  +1.09%  '

$ws.Cells.Item(36, 5).Value = '  +1.12%  '

$ws.Cells.Item(38, 5).Value = '  -0.11%  '

$ws.Cells.Item(39, 4).Value = "'" + '4.13'
$ws.Cells.Item(39, 5).Value = '  +6.18%  '

$ws.Cells.Item(40, 4).Value = "'" + '321.58'
$ws.Cells.Item(40, 5).Value = '  +11.08%  '

$ws.Cells.Item(41, 4).Value = "'" + '38.11'
$ws.Cells.Item(41, 5).Value = '  -1.07%  '

$ws.Cells.Item(42, 5).Value = '  +2.10%  '

$ws.Cells.Item(43, 4).Value = "'" + '140.61'
$ws.Cells.Item(43, 5).Value = '  -0.07%  '

$ws.Cells.Item(44, 5).Value = '  +1.21%  '

$ws.Cells.Item(45, 4).Value = "'" + '0.0946'
$ws.Cells.Item(45, 5).Value = '  -0.56%  '

$ws.Cells.Item(46, 4).Value = "'" + '19.47'
$ws.Cells.Item(46, 5).Value = '  +7.52%  '

$ws.Cells.Item(47, 4).Value = "'" + '0.0498'
$ws.Cells.Item(47, 5).Value = '  +0.67%  '

$ws.Cells.Item(48, 4).Value = "'" + '0.561'
$ws.Cells.Item(48, 5).Value = '  +0.58%  '

$ws.Cells.Item(49, 5).Value = '  +1.13%  '

$ws.Cells.Item(50, 4).Value = '0.0₆0211'
$ws.Cells.Item(50, 5).Value = '  +12.06%  '

$ws.Cells.Item(51, 5).Value = '  +0.72%  '
